$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (row 274-275),
# pushing the existing history (rows 274:409) down to (276:411).
$ws.Range("A274:A275").EntireRow.Insert()

# Row 274: new "Primera" quality record for the latest week (date 44704).
$ws.Range("A274").Value = 8
$ws.Range("B274").Value = "Terminal La Palmera de La Serena"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44704
$ws.Range("E274").Value = 4
$ws.Range("F274").Value = 100112009
$ws.Range("G274").Value = "Acelga"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 2540
$ws.Range("K274").Value = 600
$ws.Range("L274").Value = 700
$ws.Range("M274").Value = 650
$ws.Range("N274").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O274").Value = "Provincia del Elquí"
$ws.Range("P274").Value = 325
$ws.Range("Q274").Value = 2
$ws.Range("R274").Value = "Hortaliza"

# Row 275: new "Segunda" quality record for the latest week (date 44704).
$ws.Range("A275").Value = 8
$ws.Range("B275").Value = "Terminal La Palmera de La Serena"
$ws.Range("C275").Value = "Coquimbo"
$ws.Range("D275").Value = 44704
$ws.Range("E275").Value = 4
$ws.Range("F275").Value = 100112009
$ws.Range("G275").Value = "Acelga"
$ws.Range("H275").Value = "Sin especificar"
$ws.Range("I275").Value = "Segunda"
$ws.Range("J275").Value = 1400
$ws.Range("K275").Value = 500
$ws.Range("L275").Value = 550
$ws.Range("M275").Value = 525
$ws.Range("N275").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O275").Value = "Provincia del Elquí"
$ws.Range("P275").Value = 262
$ws.Range("Q275").Value = 2
$ws.Range("R275").Value = "Hortaliza"
